$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: full name built from A & B via formula, Column D: the same full
# name pasted in as a literal value (as if copy/pasted as values).
for ($r = 1; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Formula = '=A' + $r + '&" "&B' + $r
}

$fullNames = @(
    "Manuel Neuer",
    "Kylian Mbappe",
    "Lionel Messi",
    "Philipe Coutinho",
    "Zlatan Ibrahimovic",
    "Virgil Van Dijk",
    "Allison Becker",
    "Marc Van der Stegen",
    "Jordi Alba",
    "Sergio Ramos",
    "N'golo Kante",
    "Bruno Fernandes",
    "Luka Modric",
    "Ilkay Gundogan",
    "Bernardo Silva",
    "Raheem Sterling"
)

for ($r = 1; $r -le 16; $r++) {
    $ws.Cells.Item($r, 4).Value = $fullNames[$r - 1]
}

# Threaded comment on D1, matching the one already present on A1.
$ws.Range("D1").AddCommentThreaded("Lim inn verdier her") | Out-Null

# Column C now holds the widest text in the sheet ("Marc Van der Stegen"),
# so Excel auto-sizes it, same as it auto-sized column A originally.
$ws.Columns.Item(3).AutoFit() | Out-Null

# Selection / view bookkeeping to mirror the authored workbook.
$ws.Range("D10").Select()
